# Implements: insert a "Developer" column (B) into the "Detailed" sheet,
# populating it with developer names and matching the header/body formatting
# used elsewhere in the table (thin box borders, light-blue fill, centered
# wrapped text), then leaves the selection on B4 as in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Detailed")

# ---------------------------------------------------------------------
# 1. Insert a new column before column B (old B..G become C..H)
# ---------------------------------------------------------------------
$ws.Columns.Item(2).Insert() | Out-Null

# Give the new column roughly the same width as column A
$ws.Columns.Item(2).ColumnWidth = 12.33

# ---------------------------------------------------------------------
# 2. Populate the new column
# ---------------------------------------------------------------------
$ws.Range("B3").Value = "Developer"
$ws.Range("B4").Value = "Chetan"
$ws.Range("B5").Value = "Chetan"
$ws.Range("B6").Value = "Chetan"
$ws.Range("B7").Value = "Chetan"
$ws.Range("B8").Value = "Chetan"
$ws.Range("B9").Value = "Maitri"
# B10:B13 and B14:B17 are left blank (matching the original workbook)

# ---------------------------------------------------------------------
# 3. Formatting helpers
# ---------------------------------------------------------------------
function Format-BoxCell($rng) {
    $rng.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $rng.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $rng.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $rng.Borders.Item(10).LineStyle = 1  # xlEdgeRight
}

function Format-SidesCell($rng, [bool]$top, [bool]$bottom) {
    $rng.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $rng.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    if ($top) { $rng.Borders.Item(8).LineStyle = 1 }       # xlEdgeTop
    if ($bottom) { $rng.Borders.Item(9).LineStyle = 1 }    # xlEdgeBottom
}

# Light blue fill (theme Accent1, Lighter 40%) used for the new column's body
$devFillColor = 14461583   # RGB(143,170,220) == Excel's "Blue, Accent 1, Lighter 40%"
# Header fill (grey) matching the other header cells
$headerFillColor = $ws.Range("A3").Interior.Color

# ---------------------------------------------------------------------
# 4. Header cell (B3) - same look as the rest of the header row
# ---------------------------------------------------------------------
$hdr = $ws.Range("B3")
Format-BoxCell $hdr
$hdr.Interior.Color = $headerFillColor
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.WrapText = $true

# ---------------------------------------------------------------------
# 5. Body cells (B4:B17) - light blue fill, centered + wrapped text
# ---------------------------------------------------------------------
$bodyRange = $ws.Range("B4:B17")
$bodyRange.Interior.Color = $devFillColor
$bodyRange.HorizontalAlignment = -4108   # xlCenter
$bodyRange.WrapText = $true

# B4:B8 -> single outlined block per row (full box), vertically centered
foreach ($r in 4..8) {
    $c = $ws.Range("B$r")
    Format-BoxCell $c
    $c.VerticalAlignment = -4108   # xlCenter
}

# B9 -> top border only (sits below the header, no vertical centering)
$c = $ws.Range("B9")
Format-SidesCell $c $true $false

# B10:B13 -> visually grouped like A10:A13 (top/middle/middle/bottom), vcentered
$c = $ws.Range("B10"); Format-SidesCell $c $true $false;  $c.VerticalAlignment = -4108
$c = $ws.Range("B11"); Format-SidesCell $c $false $false; $c.VerticalAlignment = -4108
$c = $ws.Range("B12"); Format-SidesCell $c $false $false; $c.VerticalAlignment = -4108
$c = $ws.Range("B13"); Format-SidesCell $c $false $true;  $c.VerticalAlignment = -4108

# B14:B17 -> full box per row, vertically centered
foreach ($r in 14..17) {
    $c = $ws.Range("B$r")
    Format-BoxCell $c
    $c.VerticalAlignment = -4108   # xlCenter
}

# ---------------------------------------------------------------------
# 6. Final selection / view state
# ---------------------------------------------------------------------
$ws.Range("B4").Select() | Out-Null

$wb.Save()
